# Data entry for FvR (forbesi vs rubens) experiment, Days 3-5, Af/Ar groups,
# trials 0-8 (Af) / 0-8 (Ar) — "enter thru AvF 5"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 40,8
$data[0,0] = 'FvR'
$data[0,1] = 'Af'
$data[0,2] = 3
$data[0,3] = 0
$data[0,4] = 5
$data[0,5] = $false
$data[0,6] = $true
$data[0,7] = ""
$data[1,0] = 'FvR'
$data[1,1] = 'Af'
$data[1,2] = 3
$data[1,3] = 1
$data[1,4] = 5
$data[1,5] = $false
$data[1,6] = $true
$data[1,7] = ""
$data[2,0] = 'FvR'
$data[2,1] = 'Af'
$data[2,2] = 3
$data[2,3] = 2
$data[2,4] = 5
$data[2,5] = $false
$data[2,6] = $true
$data[2,7] = ""
$data[3,0] = 'FvR'
$data[3,1] = 'Af'
$data[3,2] = 3
$data[3,3] = 3
$data[3,4] = 4
$data[3,5] = $false
$data[3,6] = $true
$data[3,7] = '"One big glowing guy. Potential cannibal?"'
$data[4,0] = 'FvR'
$data[4,1] = 'Af'
$data[4,2] = 3
$data[4,3] = 4
$data[4,4] = 4
$data[4,5] = $false
$data[4,6] = $true
$data[4,7] = ""
$data[5,0] = 'FvR'
$data[5,1] = 'Af'
$data[5,2] = 3
$data[5,3] = 5
$data[5,4] = 3
$data[5,5] = $false
$data[5,6] = $true
$data[5,7] = ""
$data[6,0] = 'FvR'
$data[6,1] = 'Af'
$data[6,2] = 3
$data[6,3] = 6
$data[6,4] = 3
$data[6,5] = $false
$data[6,6] = $true
$data[6,7] = ""
$data[7,0] = 'FvR'
$data[7,1] = 'Af'
$data[7,2] = 3
$data[7,3] = 7
$data[7,4] = 3
$data[7,5] = $false
$data[7,6] = $true
$data[7,7] = ""
$data[8,0] = 'FvR'
$data[8,1] = 'Af'
$data[8,2] = 3
$data[8,3] = 8
$data[8,4] = 3
$data[8,5] = $true
$data[8,6] = $true
$data[8,7] = ""
$data[9,0] = 'FvR'
$data[9,1] = 'Ar'
$data[9,2] = 3
$data[9,3] = 0
$data[9,4] = 5
$data[9,5] = $false
$data[9,6] = $false
$data[9,7] = ""
$data[10,0] = 'FvR'
$data[10,1] = 'Ar'
$data[10,2] = 3
$data[10,3] = 1
$data[10,4] = 5
$data[10,5] = $false
$data[10,6] = $false
$data[10,7] = ""
$data[11,0] = 'FvR'
$data[11,1] = 'Ar'
$data[11,2] = 3
$data[11,3] = 2
$data[11,4] = 4
$data[11,5] = $false
$data[11,6] = $false
$data[11,7] = ""
$data[12,0] = 'FvR'
$data[12,1] = 'Ar'
$data[12,2] = 3
$data[12,3] = 3
$data[12,4] = 4
$data[12,5] = $false
$data[12,6] = $false
$data[12,7] = ""
$data[13,0] = 'FvR'
$data[13,1] = 'Ar'
$data[13,2] = 3
$data[13,3] = 4
$data[13,4] = 3
$data[13,5] = $false
$data[13,6] = $false
$data[13,7] = ""
$data[14,0] = 'FvR'
$data[14,1] = 'Ar'
$data[14,2] = 3
$data[14,3] = 5
$data[14,4] = 3
$data[14,5] = $false
$data[14,6] = $false
$data[14,7] = ""
$data[15,0] = 'FvR'
$data[15,1] = 'Ar'
$data[15,2] = 3
$data[15,3] = 6
$data[15,4] = 3
$data[15,5] = $false
$data[15,6] = $false
$data[15,7] = ""
$data[16,0] = 'FvR'
$data[16,1] = 'Ar'
$data[16,2] = 3
$data[16,3] = 7
$data[16,4] = 3
$data[16,5] = $false
$data[16,6] = $false
$data[16,7] = ""
$data[17,0] = 'FvR'
$data[17,1] = 'Ar'
$data[17,2] = 3
$data[17,3] = 8
$data[17,4] = 3
$data[17,5] = $true
$data[17,6] = $false
$data[17,7] = ""
$data[18,0] = 'FvR'
$data[18,1] = 'Af'
$data[18,2] = 4
$data[18,3] = 0
$data[18,4] = 5
$data[18,5] = $false
$data[18,6] = $true
$data[18,7] = ""
$data[19,0] = 'FvR'
$data[19,1] = 'Af'
$data[19,2] = 4
$data[19,3] = 1
$data[19,4] = 5
$data[19,5] = $false
$data[19,6] = $true
$data[19,7] = ""
$data[20,0] = 'FvR'
$data[20,1] = 'Af'
$data[20,2] = 4
$data[20,3] = 2
$data[20,4] = 5
$data[20,5] = $false
$data[20,6] = $true
$data[20,7] = ""
$data[21,0] = 'FvR'
$data[21,1] = 'Af'
$data[21,2] = 4
$data[21,3] = 3
$data[21,4] = 5
$data[21,5] = $false
$data[21,6] = $true
$data[21,7] = ""
$data[22,0] = 'FvR'
$data[22,1] = 'Af'
$data[22,2] = 4
$data[22,3] = 4
$data[22,4] = 5
$data[22,5] = $false
$data[22,6] = $true
$data[22,7] = ""
$data[23,0] = 'FvR'
$data[23,1] = 'Af'
$data[23,2] = 4
$data[23,3] = 5
$data[23,4] = 5
$data[23,5] = $true
$data[23,6] = $true
$data[23,7] = ""
$data[24,0] = 'FvR'
$data[24,1] = 'Ar'
$data[24,2] = 4
$data[24,3] = 0
$data[24,4] = 5
$data[24,5] = $false
$data[24,6] = $false
$data[24,7] = ""
$data[25,0] = 'FvR'
$data[25,1] = 'Ar'
$data[25,2] = 4
$data[25,3] = 1
$data[25,4] = 2
$data[25,5] = $false
$data[25,6] = $false
$data[25,7] = '"Non-glowing ossicles spotted"'
$data[26,0] = 'FvR'
$data[26,1] = 'Ar'
$data[26,2] = 4
$data[26,3] = 2
$data[26,4] = 0
$data[26,5] = $false
$data[26,6] = $false
$data[26,7] = ""
$data[27,0] = 'FvR'
$data[27,1] = 'Ar'
$data[27,2] = 4
$data[27,3] = 3
$data[27,4] = 0
$data[27,5] = $false
$data[27,6] = $false
$data[27,7] = ""
$data[28,0] = 'FvR'
$data[28,1] = 'Ar'
$data[28,2] = 4
$data[28,3] = 4
$data[28,4] = 0
$data[28,5] = $false
$data[28,6] = $false
$data[28,7] = ""
$data[29,0] = 'FvR'
$data[29,1] = 'Ar'
$data[29,2] = 4
$data[29,3] = 5
$data[29,4] = 0
$data[29,5] = $true
$data[29,6] = $false
$data[29,7] = ""
$data[30,0] = 'FvR'
$data[30,1] = 'Af'
$data[30,2] = 5
$data[30,3] = 0
$data[30,4] = 5
$data[30,5] = $false
$data[30,6] = $true
$data[30,7] = ""
$data[31,0] = 'FvR'
$data[31,1] = 'Af'
$data[31,2] = 5
$data[31,3] = 1
$data[31,4] = 5
$data[31,5] = $false
$data[31,6] = $true
$data[31,7] = ""
$data[32,0] = 'FvR'
$data[32,1] = 'Af'
$data[32,2] = 5
$data[32,3] = 2
$data[32,4] = 4
$data[32,5] = $false
$data[32,6] = $true
$data[32,7] = ""
$data[33,0] = 'FvR'
$data[33,1] = 'Af'
$data[33,2] = 5
$data[33,3] = 3
$data[33,4] = 2
$data[33,5] = $false
$data[33,6] = $true
$data[33,7] = ""
$data[34,0] = 'FvR'
$data[34,1] = 'Af'
$data[34,2] = 5
$data[34,3] = 4
$data[34,4] = 1
$data[34,5] = $true
$data[34,6] = $true
$data[34,7] = ""
$data[35,0] = 'FvR'
$data[35,1] = 'Ar'
$data[35,2] = 5
$data[35,3] = 0
$data[35,4] = 5
$data[35,5] = $false
$data[35,6] = $false
$data[35,7] = ""
$data[36,0] = 'FvR'
$data[36,1] = 'Ar'
$data[36,2] = 5
$data[36,3] = 1
$data[36,4] = 5
$data[36,5] = $false
$data[36,6] = $false
$data[36,7] = ""
$data[37,0] = 'FvR'
$data[37,1] = 'Ar'
$data[37,2] = 5
$data[37,3] = 2
$data[37,4] = 2
$data[37,5] = $false
$data[37,6] = $false
$data[37,7] = ""
$data[38,0] = 'FvR'
$data[38,1] = 'Ar'
$data[38,2] = 5
$data[38,3] = 3
$data[38,4] = 1
$data[38,5] = $false
$data[38,6] = $false
$data[38,7] = ""
$data[39,0] = 'FvR'
$data[39,1] = 'Ar'
$data[39,2] = 5
$data[39,3] = 4
$data[39,4] = 1
$data[39,5] = $true
$data[39,6] = $false
$data[39,7] = ""

$ws.Range("A42:H81").Value = $data

# Row 78 note contains the italic species name "rubens" embedded in the sentence.
$noteCell = $ws.Range("H78")
$noteCell.Value = '"Non-glower (rubens) being cannibalized"'
$italicRun = $noteCell.Characters(14, 6)
$italicRun.Font.Italic = $true
$tailRun = $noteCell.Characters(20, 21)
$tailRun.Font.Name = "Aptos Narrow"
$tailRun.Font.Size = 12

# Restore view state to match where data entry left off.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$win.ScrollColumn = 1
$ws.Range("B86").Select()
